# LogicComponentSequenceDiagram.pptx edit
#
# The authoritative diff also touches the "datetimeFigureOut" footer
# field on many slides (3/23/18 -> 4/12/18), but this deck only carries
# a single slide and that slide has no date/footer placeholder at all,
# so there is nothing to change there.
#
# The remaining hunks all live on slide 1 and are simple "typo fix"
# style edits where two or three adjacent <a:r> runs collapse into one
# run (PowerPoint merges runs automatically when you select the exact
# text span covering them and retype it). We reproduce that by
# addressing the text with TextRange.Characters(start,length) and
# assigning .Text - when the span exactly covers a set of whole runs,
# the engine consolidates them into a single run using the formatting
# of the run that "owns" the selection start. For the handful of spots
# where the run carrying the desired (error-free) formatting is not
# the first one, we first delete the earlier run(s)' characters, use
# InsertBefore() on the remaining run (which inherits its formatting),
# and then re-merge the whole range.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$ldquo = [char]0x201C  # "  U+201C LEFT DOUBLE QUOTATION MARK

# --- Shape 2 ("Rectangle 62"): ":" + "LogicManager" -> ":LogicManager"
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = ":LogicManager"

# --- Shape 5 ("Rectangle 62"): first paragraph ":Address" -> ":Internship"
$sh = $s.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(1, 8)
$c.Text = ":Internship"

# --- Shape 8 ("Rectangle 62"): "s" + ":Save" -> "s:Save" (before the line break + "Command")
$sh = $s.Shapes.Item(8)
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(1, 6)
$c.Text = "s:Save"

# --- Shape 13 ("TextBox 25"): execute("  + save + " " -> execute("save "
$sh = $s.Shapes.Item(13)
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(8, 7)
$c.Text = "(" + $ldquo + "save "

# --- Shape 27 ("TextBox 77"): "update" + "Person" + "(p)" -> "updatePerson(p)"
# The error-free formatting lives on the 3rd run ("(p)"), so delete the
# first two runs' characters and re-insert the missing text in front of
# the remaining run (inherits its rPr), then consolidate into one run.
$sh = $s.Shapes.Item(27)
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(1, 12)
$c.Delete()
$tr = $sh.TextFrame.TextRange
$tr.InsertBefore("updatePerson")
$tr = $sh.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = "updatePerson(p)"

# --- Shape 29 ("TextBox 79"): parse + ("  + save + " " -> parse + ("save " + 1")
$sh = $s.Shapes.Item(29)
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(6, 7)
$c.Text = "(" + $ldquo + "save "

# --- Shape 32 ("Rectangle 62"): "result:Command" + " Result" -> "result:Command Result"
# Error-free formatting lives on the 2nd run (" Result"), so delete the
# first run's characters and re-insert in front of what remains.
$sh = $s.Shapes.Item(32)
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(1, 14)
$c.Delete()
$tr = $sh.TextFrame.TextRange
$tr.InsertBefore("result:Command")
$tr = $sh.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = "result:Command Result"

# --- Shape 36 ("Rectangle 62"): ":" + "Save" + "Command" -> ":SaveCommand"
$sh = $s.Shapes.Item(36)
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(1, 12)
$c.Text = ":SaveCommand"

# --- Shape 49 ("TextBox 68"): "addSavedTagToInternship" + "()" -> "addSavedTagToInternship()"
# Error-free formatting lives on the 2nd run ("()"), so delete the
# first run's characters and re-insert in front of what remains.
$sh = $s.Shapes.Item(49)
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(1, 23)
$c.Delete()
$tr = $sh.TextFrame.TextRange
$tr.InsertBefore("addSavedTagToInternship")
$tr = $sh.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = "addSavedTagToInternship()"
